$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-05 Wednesday" "2025-11-06 Thursday"
Replace-Text "683×2=1366" "887×9=7983"
Replace-Text "669×9=6021" "336×9=3024"
Replace-Text "324×4=1296" "741×4=2964"
Replace-Text "911×3=2733" "554×8=4432"
Replace-Text "423×4=1692" "241×7=1687"
Replace-Text "564×8=4512" "852×2=1704"
Replace-Text "269×3=807" "162×2=324"
Replace-Text "389×3=1167" "819×5=4095"
Replace-Text "968×4=3872" "205×5=1025"
Replace-Text "825×6=4950" "746×6=4476"
Replace-Text "302×6=1812" "858×7=6006"
Replace-Text "441×2=882" "367×4=1468"
Replace-Text "501×7=3507" "762×5=3810"
Replace-Text "625×2=1250" "677×7=4739"
Replace-Text "475×5=2375" "473×7=3311"
Replace-Text "844×9=7596" "929×5=4645"
Replace-Text "878×7=6146" "356×8=2848"
Replace-Text "321×8=2568" "315×9=2835"
Replace-Text "137×5=685" "283×2=566"
Replace-Text "415×7=2905" "881×6=5286"
Replace-Text "742×4=2968" "929×7=6503"
Replace-Text "675×4=2700" "701×9=6309"
Replace-Text "498×9=4482" "754×3=2262"
Replace-Text "115×5=575" "673×8=5384"
Replace-Text "656×5=3280" "220×3=660"
